# "Actualizacion plan de calidad"
# The row that tracked the "Requerimientos" deliverable (B19:D19 - previously
# "Requerimientos" / "Al finalizar etapa de ventas" / "Jovanny Zepeda") is
# cleared out, leaving the row blank but keeping its existing formatting.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Plan de Calidad")

$ws.Range("B19:D19").ClearContents()

# Leave the workbook scrolled/selected where the author left off editing.
$ws.Activate()
$win = $excel.ActiveWindow
$win.ScrollRow = 6
$win.ScrollColumn = 1
$ws.Range("E19").Select()
